# Apply updated monitoring data (huambocancha alta) - new totals and re-ranked order
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering (by row) with name and total_registros, reflecting the
# re-ranked / updated values from the diff.
$data = @(
    @{ Name = "INCIO SANCHEZ PAOLA KATHERINE";       Total = 101 },
    @{ Name = "GUEVARA IDROGO DENNIS PERCY";          Total = 97 },
    @{ Name = "TANTALEAN BUSTAMANTE ESTALIN YOEL";    Total = 95 },
    @{ Name = "HUAYHUA VALDIVIA LUZ EXMILDA";         Total = 86 },
    @{ Name = "LINARES PEREZ YANASELY";               Total = 86 },
    @{ Name = "PEREZ LINARES TATHIANA";               Total = 86 },
    @{ Name = "MONDRAGON HERNANDEZ WILMER JUNIOR";    Total = 85 },
    @{ Name = "MEDINA TAPIA ANA YULI";                Total = 85 },
    @{ Name = "CAMPOS PEREZ YOVERLY";                 Total = 85 },
    @{ Name = "CHAVEZ VILLANUEVA SILVIA JANETH";      Total = 84 },
    @{ Name = "DELGADO VASQUEZ FLOR MAGALY";          Total = 84 },
    @{ Name = "LOZADA ROJAS LUZ ELENA";               Total = 84 },
    @{ Name = "SOTO LOZANO LUZDINA";                  Total = 79 },
    @{ Name = "VASQUEZ SILVA ALOIS ADOLF";             Total = 73 }
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i].Name
    $ws.Cells.Item($row, 2).Value = $data[$i].Total
}
